$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (full-row content swap/permutation)
$rowMap = @{2=6; 3=2; 6=3; 7=13; 8=14; 9=10; 10=7; 12=8; 13=9; 14=12; 15=19; 16=20; 17=15; 18=16; 19=17; 20=18; 28=32; 30=28; 31=30; 32=31; 40=48; 41=40; 42=43; 43=42; 44=49; 45=44; 46=45; 47=41; 48=46; 49=47; 50=52; 51=53; 52=54; 53=51; 54=50; 56=58; 57=56; 58=57; 70=71; 71=70; 72=73; 73=74; 74=72; 75=76; 76=75; 77=78; 78=77; 81=85; 82=83; 83=81; 84=82; 85=84; 90=92; 91=94; 92=90; 93=91; 94=93; 97=101; 98=102; 99=97; 100=98; 101=103; 102=100; 103=99}

# Snapshot original row contents BEFORE any writes, split to avoid Y/AA (date-looking text)
# getting reinterpreted as real dates by Excel when written back via COM.
$snapA = @{}   # columns A:X  (1-24)
$snapZ = @{}   # column  Z    (26)
$snapAB = @{}  # columns AB:AY (28-51)
$involvedRows = @(2; 3; 6; 7; 8; 9; 10; 12; 13; 14; 15; 16; 17; 18; 19; 20; 28; 30; 31; 32; 40; 41; 42; 43; 44; 45; 46; 47; 48; 49; 50; 51; 52; 53; 54; 56; 57; 58; 70; 71; 72; 73; 74; 75; 76; 77; 78; 81; 82; 83; 84; 85; 90; 91; 92; 93; 94; 97; 98; 99; 100; 101; 102; 103)
foreach ($r in $involvedRows) {
    $snapA[$r] = $ws.Range("A" + $r + ":X" + $r).Value2
    $snapZ[$r] = $ws.Range("Z" + $r + ":Z" + $r).Value2
    $snapAB[$r] = $ws.Range("AB" + $r + ":AY" + $r).Value2
}

# Write back: each destination row receives the snapshot of its source row
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $ws.Range("A" + $destRow + ":X" + $destRow).Value2 = $snapA[$srcRow]
    $ws.Range("Z" + $destRow + ":Z" + $destRow).Value2 = $snapZ[$srcRow]
    $ws.Range("AB" + $destRow + ":AY" + $destRow).Value2 = $snapAB[$srcRow]
}

Write-Host "Done."
